# Updates "grouped matches" sheet (BDRC ID / 84000 ID matches) rows 2-60.
# New data reflects the 84000 ids added and corrected BDRC<->84000 groupings.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grouped matches")

$data = @(
    @(2, 'P4CZ16780 ', '{''eft:manjusrigarbha''}'),
    @(3, 'P8213 ', '{''eft:t-vidyakarasimha''}'),
    @(4, 'P0TMP092', '{''eft:anandasri-s-''}'),
    @(5, 'P8219', '{''eft:visuddhasimha''}'),
    @(6, 'P8268', '{''eft:buddhaprabha''}'),
    @(7, 'P8209', '{''eft:dzi-na-mi-tra-k-'', ''eft:jinamitra'', ''eft:jinamitra-k-''}'),
    @(8, 'P8205', '{''eft:yesh-d-ye-shes-sde-'', ''eft:yesh-d-'', ''eft:band-yesh-de'', ''eft:ye-shes-sde'', ''eft:zhang-yesh-d-''}'),
    @(9, 'https://lod.dila.edu.tw/resource.php?id=A000089', '{''eft:siladharma''}'),
    @(10, 'P00KG07267', '{''eft:sarvajnadeva'', ''eft:sarvanyadeva''}'),
    @(11, 'P4259', '{''eft:dpal-gyi-lhun-po'', ''eft:ban-de-dpal-gyi-lhun-po'', ''eft:palgyi-lh-npo''}'),
    @(12, 'P3285 ', '{''eft:sakya-yesh-''}'),
    @(13, 'P753', '{''eft:rin-chen-bzang-po''}'),
    @(14, 'P4CZ15137', '{''eft:kumarakalasa''}'),
    @(15, 'P8267', '{''eft:vijayasila''}'),
    @(16, 'P0TMP098', '{''eft:jinavara''}'),
    @(17, 'P0TMP080', '{''eft:hwa-shang-zab-mo''}'),
    @(18, 'P8210', '{''eft:danasila''}'),
    @(19, 'P3709 ', '{''eft:phakpa-sherab''}'),
    @(20, 'P8269', '{''eft:dgon-gling-rma''}'),
    @(21, 'P4255', '{''eft:ye-shes-snying-po'', ''eft:yesh-nyingpo'', ''eft:t-jnanagarbha''}'),
    @(22, 'P4242', '{''eft:sherab-lekpa''}'),
    @(23, 'P8222', '{''eft:jnanasiddhi''}'),
    @(24, 'P0RK8', '{''eft:dharmapala''}'),
    @(25, 'P8151', '{''eft:gayadhara''}'),
    @(26, 'P8222 ', '{''eft:jnanasidhi''}'),
    @(27, 'P8217', '{''eft:t-jnanagarbha'', ''eft:jnanagarbha''}'),
    @(28, 'P8205 ', '{''eft:band-yesh-d-''}'),
    @(29, 'P3214 ', '{''eft:danasila''}'),
    @(30, 'P8206', '{''eft:celu''}'),
    @(31, 'P8245', '{''eft:buddhakaravarma''}'),
    @(32, 'P4263', '{''eft:dge-ba-dpal''}'),
    @(33, 'P2548', '{''eft:prajnavarma'', ''eft:prajnavarman''}'),
    @(34, 'P8261', '{''eft:munivarma'', ''eft:munivarman''}'),
    @(35, 'P8249', '{''eft:dharmakara''}'),
    @(36, 'P8093', '{''eft:kamalagupta''}'),
    @(37, '?', '{''eft:sakyasena''}'),
    @(38, 'P8263', '{''eft:leki-d-''}'),
    @(39, 'P8228', '{''eft:surendrabodhi''}'),
    @(40, 'P8220', '{''eft:devacandra''}'),
    @(41, 'P1KG8854', '{''eft:surendrabodhi'', ''eft:silendrabodhi''}'),
    @(42, 'P8171', '{''eft:dharmasribhadra''}'),
    @(43, 'P0TMPT007', '{''eft:rnam-par-mi-rtog-pa''}'),
    @(44, 'P8260', '{''eft:dpal-dbyangs''}'),
    @(45, 'P4CZ16819', '{''eft:sakyaprabha''}'),
    @(46, 'P8273', '{''eft:rinchen-tso'', ''eft:rin-chen-tsho''}'),
    @(47, 'P0TMP104', '{''eft:punyasambhava''}'),
    @(48, 'P8213', '{''eft:vidyakarasimha''}'),
    @(49, 'P5651', '{''eft:pa-tshab-nyi-ma-grags''}'),
    @(50, 'P8183', '{''eft:klu-i-rgyal-mtshan'', ''eft:cog-ro-klu-i-rgyal-mtshan''}'),
    @(51, 'P4258', '{''eft:dpal-byor''}'),
    @(52, 'P8211', '{''eft:vidyakaraprabha''}'),
    @(53, 'P8182', '{''eft:paltsek'', ''eft:ska-ba-dpal-brtsegs'', ''eft:ban-de-dpal-brtsegs'', ''eft:kawa-paltsek-under-the-name-paltsek-raksita-'', ''eft:dpal-brtsegs''}'),
    @(54, 'P8265', '{''eft:ratnaraksita''}'),
    @(55, 'P2956', '{''eft:krsnapandita''}'),
    @(56, 'P1KG8854 ', '{''eft:srilendrabodhi''}'),
    @(57, 'P2637', '{''eft:trakpa-gyaltsen''}'),
    @(58, 'P3379', '{''eft:dipamkarasrijnana'', ''eft:dipamkara-srijnana''}'),
    @(59, 'P8266', '{''eft:dharmatasila'', ''eft:ch-nyi-tsultrim''}'),
    @(60, 'P3456', '{''eft:tshul-khrims-rgyal-ba''}')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
